# Update cryptocurrency price/volume data to reflect latest market snapshot
# (and fix the swapped FirstDigitalUSD / ImmutableX row ordering).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.851.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.22%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.982.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.37%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'541.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.25%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'136.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.32%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'2.978.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.36%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.485"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.44%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +10.60%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.17%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.442"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.86%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000218"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.36%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'33.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.09%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.449.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.78%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'61.792.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.42%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.36%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.980.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.42%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.51%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'465.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.12%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.20%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.649"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.89%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'79.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.70%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'12.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.10%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.01%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.57%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.78%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'ImmutableX"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.44%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'FirstDigitalUSD"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.21%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'25.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.26%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.15%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'2.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.33%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.21%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'54.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.60%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.75%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'449.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.13%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0801"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.43%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0384"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'2.936.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -9.13%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.29%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'7.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.93%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'26.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.20%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D46").Value = "'0.245"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.12%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.40%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.32%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'114.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.40%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0₃0487"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.05%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -2.35%  "
$ws.Range("E51").Style = "Normal"
